{"js": "// Replace the worksheet date and each \"A\u00f7B=C, D\" answer cell with the new\n// values from the commit. Every old value is unique in the document, so a\n// literal, case-sensitive body.search() + insertText(\"Replace\") round trip\n// is sufficient and keeps all existing run formatting (font/size) intact.\nconst replacements = [\n  [\"2025-06-27 Friday\", \"2025-06-28 Saturday\"],\n  [\"126\u00f74=31, 2\", \"630\u00f75=126, 0\"],\n  [\"707\u00f77=101, 0\", \"988\u00f74=247, 0\"],\n  [\"804\u00f75=160, 4\", \"299\u00f79=33, 2\"],\n  [\"958\u00f76=159, 4\", \"114\u00f74=28, 2\"],\n  [\"913\u00f72=456, 1\", \"677\u00f76=112, 5\"],\n  [\"805\u00f79=89, 4\", \"268\u00f76=44, 4\"],\n  [\"995\u00f79=110, 5\", \"573\u00f78=71, 5\"],\n  [\"326\u00f79=36, 2\", \"931\u00f77=133, 0\"],\n  [\"116\u00f75=23, 1\", \"943\u00f77=134, 5\"],\n  [\"895\u00f73=298, 1\", \"399\u00f79=44, 3\"],\n  [\"256\u00f77=36, 4\", \"108\u00f75=21, 3\"],\n  [\"677\u00f79=75, 2\", \"131\u00f74=32, 3\"],\n  [\"815\u00f79=90, 5\", \"120\u00f76=20, 0\"],\n  [\"274\u00f73=91, 1\", \"195\u00f78=24, 3\"],\n  [\"603\u00f78=75, 3\", \"890\u00f79=98, 8\"],\n  [\"830\u00f75=166, 0\", \"693\u00f72=346, 1\"],\n  [\"736\u00f73=245, 1\", \"630\u00f74=157, 2\"],\n  [\"153\u00f77=21, 6\", \"113\u00f77=16, 1\"],\n  [\"333\u00f78=41, 5\", \"610\u00f73=203, 1\"],\n  [\"444\u00f76=74, 0\", \"696\u00f77=99, 3\"],\n  [\"626\u00f72=313, 0\", \"465\u00f73=155, 0\"],\n  [\"999\u00f78=124, 7\", \"456\u00f78=57, 0\"],\n  [\"425\u00f76=70, 5\", \"765\u00f73=255, 0\"],\n  [\"807\u00f77=115, 2\", \"827\u00f78=103, 3\"],\n  [\"289\u00f79=32, 1\", \"911\u00f77=130, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and each \"A\u00f7B=C, D\" answer cell with the new\n# values from the commit. Every old value is unique in the document, so a\n# literal, case-sensitive Find/Replace round trip (wdReplaceAll, though each\n# only ever matches once) is sufficient and leaves existing run formatting\n# (font/size) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-06-27 Friday\", \"2025-06-28 Saturday\"),\n    @(\"126\u00f74=31, 2\", \"630\u00f75=126, 0\"),\n    @(\"707\u00f77=101, 0\", \"988\u00f74=247, 0\"),\n    @(\"804\u00f75=160, 4\", \"299\u00f79=33, 2\"),\n    @(\"958\u00f76=159, 4\", \"114\u00f74=28, 2\"),\n    @(\"913\u00f72=456, 1\", \"677\u00f76=112, 5\"),\n    @(\"805\u00f79=89, 4\", \"268\u00f76=44, 4\"),\n    @(\"995\u00f79=110, 5\", \"573\u00f78=71, 5\"),\n    @(\"326\u00f79=36, 2\", \"931\u00f77=133, 0\"),\n    @(\"116\u00f75=23, 1\", \"943\u00f77=134, 5\"),\n    @(\"895\u00f73=298, 1\", \"399\u00f79=44, 3\"),\n    @(\"256\u00f77=36, 4\", \"108\u00f75=21, 3\"),\n    @(\"677\u00f79=75, 2\", \"131\u00f74=32, 3\"),\n    @(\"815\u00f79=90, 5\", \"120\u00f76=20, 0\"),\n    @(\"274\u00f73=91, 1\", \"195\u00f78=24, 3\"),\n    @(\"603\u00f78=75, 3\", \"890\u00f79=98, 8\"),\n    @(\"830\u00f75=166, 0\", \"693\u00f72=346, 1\"),\n    @(\"736\u00f73=245, 1\", \"630\u00f74=157, 2\"),\n    @(\"153\u00f77=21, 6\", \"113\u00f77=16, 1\"),\n    @(\"333\u00f78=41, 5\", \"610\u00f73=203, 1\"),\n    @(\"444\u00f76=74, 0\", \"696\u00f77=99, 3\"),\n    @(\"626\u00f72=313, 0\", \"465\u00f73=155, 0\"),\n    @(\"999\u00f78=124, 7\", \"456\u00f78=57, 0\"),\n    @(\"425\u00f76=70, 5\", \"765\u00f73=255, 0\"),\n    @(\"807\u00f77=115, 2\", \"827\u00f78=103, 3\"),\n    @(\"289\u00f79=32, 1\", \"911\u00f77=130, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
